$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Docentes responsaveis:" value row (row 13), which shifts
# every subsequent row up by one and drops the final (bibliography body) row,
# matching the new A1:C21 dimension.
$ws.Rows("13:13").Delete()

# After the shift, a handful of value cells (column B/C) need to be
# re-pointed to the text that now actually belongs in those rows.
$ws.Cells.Item(10, 2).Value2 = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Cells.Item(10, 3).Value2 = "1720367 - Teresa Cristina Brazil de Paiva"

$ws.Cells.Item(13, 2).Value2 = "Semestral"
$ws.Cells.Item(13, 3).Value2 = "Semestral"

$ws.Cells.Item(15, 2).Value2 = "01/01/2020"
$ws.Cells.Item(15, 3).Value2 = "01/01/2020"

$ws.Cells.Item(18, 2).Value2 = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Cells.Item(18, 3).Value2 = "1720367 - Teresa Cristina Brazil de Paiva"

$ws.Cells.Item(19, 2).Value2 = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."
$ws.Cells.Item(19, 3).Value2 = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."

$ws.Cells.Item(20, 2).Value2 = "Média ponderada das notas atribuídas à prova, exercício e relatório."
$ws.Cells.Item(20, 3).Value2 = "Média ponderada das notas atribuídas à prova, exercício e relatório."

$ws.Cells.Item(21, 2).Value2 = "Nota final: NF ≥ 5,0"
$ws.Cells.Item(21, 3).Value2 = "Nota final: NF ≥ 5,0"
